# Revise login page ui elements
# Appends one new log row (row 17) to each of the three "...-LIFTER" sheets
# whose data currently stops at row 16 (ROW50-FE-LIFTER, ROW11-FE-LIFTER,
# ROW11-MID-LIFTER). ROW50-MID-LIFTER already extends past row 17 and is
# left untouched.

$wb = $excel.ActiveWorkbook

function Add-LogRow {
    param($ws, $row, $timeValue, $colB, $colC, $colD, $colE, $colF, $colG, $colH, $colI)

    $prevRow = $row - 1

    # Match column A's existing date/time number format so the new cell
    # lands on the same style index as the rest of the column.
    $ws.Cells.Item($row, 1).NumberFormat = $ws.Cells.Item($prevRow, 1).NumberFormat

    # Assign as numeric-looking strings so they parse into numbers with
    # full precision instead of relying on a literal with an exponent
    # (the script host's expression parser rejects `e`/`E`
    # scientific-notation number literals).
    $ws.Cells.Item($row, 1).Value = $timeValue
    $ws.Cells.Item($row, 2).Value = $colB
    $ws.Cells.Item($row, 3).Value = $colC
    $ws.Cells.Item($row, 4).Value = $colD
    $ws.Cells.Item($row, 5).Value = $colE
    $ws.Cells.Item($row, 6).Value = $colF
    $ws.Cells.Item($row, 7).Value = $colG
    $ws.Cells.Item($row, 8).Value = $colH
    $ws.Cells.Item($row, 9).Value = $colI
}

# Column G's huge ID_DEC value round-trips to the same double either way,
# but writing it as a plain (non-exponent-looking) digit string keeps the
# cell on the default "General" style instead of Excel auto-applying a
# scientific-notation number format (and allocating a brand new style) the
# way it would for a string that already looks like "5.68...e+23".
$wsFe50 = $wb.Worksheets.Item("ROW50-FE-LIFTER")
Add-LogRow $wsFe50 17 "45732.61619851852" "0x01,0x90" "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x15,0x41,0x0c," "0x01,0x86" "0xe" "400" "568631262647114000000000" "390" "14"

$wsFe11 = $wb.Worksheets.Item("ROW11-FE-LIFTER")
Add-LogRow $wsFe11 17 "45732.63623451389" "0x01,0x90" "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x1a,0x41,0x0c," "0x01,0x86" "0x14" "400" "568631262647114000000000" "390" "20"

$wsMid11 = $wb.Worksheets.Item("ROW11-MID-LIFTER")
Add-LogRow $wsMid11 17 "45732.78172160879" "0x01,0x90" "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x11,0x40,0x0c," "0x01,0x86" "0x19" "400" "568631262647114000000000" "390" "25"
